$wb = $excel.ActiveWorkbook

# OFF sheet - row 2 (Home) values
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 442
$wsOff.Range("C2").Value = 314
$wsOff.Range("D2").Value = 124
$wsOff.Range("E2").Value = 62
$wsOff.Range("F2").Value = 7
$wsOff.Range("G2").Value = 5

# DEF sheet - row 2 (Home) values
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 519
$wsDef.Range("C2").Value = 350
$wsDef.Range("D2").Value = 123
$wsDef.Range("E2").Value = 45
$wsDef.Range("F2").Value = 9
$wsDef.Range("G2").Value = 12
